# Updates cryptos list values (Price and Volume(1h) columns) to latest figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '41.528.92'
$r.Style = "Normal"

$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = '  +0.04%  '
$r.Style = "Normal"

$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '2.470.22'
$r.Style = "Normal"

$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = '  -0.73%  '
$r.Style = "Normal"

$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = '  -0.36%  '
$r.Style = "Normal"

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '314.65'
$r.Style = "Normal"

$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = '  -0.13%  '
$r.Style = "Normal"

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '91.99'
$r.Style = "Normal"

$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = '  -2.63%  '
$r.Style = "Normal"

$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = '  +0.05%  '
$r.Style = "Normal"

$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = '  -0.34%  '
$r.Style = "Normal"

$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = '  +3.33%  '
$r.Style = "Normal"

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '32.56'
$r.Style = "Normal"

$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = '  -3.11%  '
$r.Style = "Normal"

$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.0792'
$r.Style = "Normal"

$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = '  +1.12%  '
$r.Style = "Normal"

$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = '  +0.44%  '
$r.Style = "Normal"

$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '2.850.25'
$r.Style = "Normal"

$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = '  -0.82%  '
$r.Style = "Normal"

$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = '  -0.98%  '
$r.Style = "Normal"

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '16.01'
$r.Style = "Normal"

$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = '  +3.53%  '
$r.Style = "Normal"

$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '2.476.88'
$r.Style = "Normal"

$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = '  -1.36%  '
$r.Style = "Normal"

$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '0.778'
$r.Style = "Normal"

$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = '  -1.83%  '
$r.Style = "Normal"

$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '41.565.73'
$r.Style = "Normal"

$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = '  +0.20%  '
$r.Style = "Normal"

$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '6.52'
$r.Style = "Normal"

$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = '  +2.47%  '
$r.Style = "Normal"

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '0.0₃0944'
$r.Style = "Normal"

$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = '  +1.89%  '
$r.Style = "Normal"

$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '71.12'
$r.Style = "Normal"

$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = '  +1.64%  '
$r.Style = "Normal"

$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '11.13'
$r.Style = "Normal"

$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = '  -1.15%  '
$r.Style = "Normal"

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '237.36'
$r.Style = "Normal"

$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = '  +0.23%  '
$r.Style = "Normal"

$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = '  -1.65%  '
$r.Style = "Normal"

$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = '  +0.69%  '
$r.Style = "Normal"

$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = '  -0.01%  '
$r.Style = "Normal"

$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '24.77'
$r.Style = "Normal"

$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = '  +2.42%  '
$r.Style = "Normal"

$r = $ws.Range("E28")
$r.NumberFormat = "@"
$r.Value = '  -1.22%  '
$r.Style = "Normal"

$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '9.72'
$r.Style = "Normal"

$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '35.52'
$r.Style = "Normal"

$r = $ws.Range("E30")
$r.NumberFormat = "@"
$r.Value = '  -4.46%  '
$r.Style = "Normal"

$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '155.79'
$r.Style = "Normal"

$r = $ws.Range("E31")
$r.NumberFormat = "@"
$r.Value = '  +0.99%  '
$r.Style = "Normal"

$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '5.46'
$r.Style = "Normal"

$r = $ws.Range("E32")
$r.NumberFormat = "@"
$r.Value = '  -1.06%  '
$r.Style = "Normal"

$r = $ws.Range("E33")
$r.NumberFormat = "@"
$r.Value = '  +0.12%  '
$r.Style = "Normal"

$r = $ws.Range("E34")
$r.NumberFormat = "@"
$r.Value = '  +0.24%  '
$r.Style = "Normal"

$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '17.31'
$r.Style = "Normal"

$r = $ws.Range("E35")
$r.NumberFormat = "@"
$r.Value = '  -3.91%  '
$r.Style = "Normal"

$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '2.36'
$r.Style = "Normal"

$r = $ws.Range("E36")
$r.NumberFormat = "@"
$r.Value = '  -2.81%  '
$r.Style = "Normal"

$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '2.90'
$r.Style = "Normal"

$r = $ws.Range("E37")
$r.NumberFormat = "@"
$r.Value = '  -5.93%  '
$r.Style = "Normal"

$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = '  +2.83%  '
$r.Style = "Normal"

$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = '  -0.28%  '
$r.Style = "Normal"

$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = '  -4.20%  '
$r.Style = "Normal"

$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '4.00'
$r.Style = "Normal"

$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = '  -3.13%  '
$r.Style = "Normal"

$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = '  -0.47%  '
$r.Style = "Normal"

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '1.946.17'
$r.Style = "Normal"

$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = '  -2.16%  '
$r.Style = "Normal"

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '18.99'
$r.Style = "Normal"

$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = '  -4.65%  '
$r.Style = "Normal"

$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = '  -1.13%  '
$r.Style = "Normal"

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '2.92'
$r.Style = "Normal"

$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = '  -3.73%  '
$r.Style = "Normal"

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '9.11'
$r.Style = "Normal"

$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = '  +2.89%  '
$r.Style = "Normal"

$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '2.709.00'
$r.Style = "Normal"

$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = '  -0.96%  '
$r.Style = "Normal"

$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = '  +0.11%  '
$r.Style = "Normal"

$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '67.33'
$r.Style = "Normal"

$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = '  -3.27%  '
$r.Style = "Normal"

$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '52.52'
$r.Style = "Normal"

$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = '  +2.97%  '
$r.Style = "Normal"

